# Fruta / hortaliza, semanal
# Insert a new data row at row 52 (shifting all subsequent rows down by one)
# and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 52. This shifts rows 52:143 down to 53:144
# and duplicates formatting from the row being pushed down (consistent with the rest
# of the table).
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new record. Columns A, B, C, E-K hold
# constant values shared by every row in this table.
$ws.Cells.Item(52, 1).Value = 6
$ws.Cells.Item(52, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(52, 3).Value = "Metropolitana"
$ws.Cells.Item(52, 4).Value = 44571
$ws.Cells.Item(52, 5).Value = 13
$ws.Cells.Item(52, 6).Value = "Fruta"
$ws.Cells.Item(52, 7).Value = 100101
$ws.Cells.Item(52, 8).Value = "Berries"
$ws.Cells.Item(52, 9).Value = 100101004
$ws.Cells.Item(52, 10).Value = "Frambuesa"
$ws.Cells.Item(52, 11).Value = "Sin especificar"
$ws.Cells.Item(52, 12).Value = "Especial"
$ws.Cells.Item(52, 13).Value = 100
$ws.Cells.Item(52, 14).Value = 8000
$ws.Cells.Item(52, 15).Value = 8000
$ws.Cells.Item(52, 16).Value = 8000
$ws.Cells.Item(52, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(52, 18).Value = "Provincia de Linares"
$ws.Cells.Item(52, 19).Value = 4000
$ws.Cells.Item(52, 20).Value = 2

# Ensure the date cell keeps the same date/time numeric format used throughout column D.
$ws.Cells.Item(52, 4).NumberFormat = $ws.Cells.Item(53, 4).NumberFormat
